$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 351.75
$ws.Range("I9").Value = 264
$ws.Range("K9").Value = 264
$ws.Range("M9").Value = -95

$ws.Range("H135").Value = 1839.7333
$ws.Range("I135").Value = 1620.5
$ws.Range("J135").Value = 2716.6667
$ws.Range("K135").Value = 14584.5
$ws.Range("L135").Value = 24450.0003
$ws.Range("M135").Value = -12049.5
$ws.Range("N135").Value = -29520.0003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3109
$ws.Range("I61").Value = 2935.55
$ws.Range("J61").Value = 3976.25
$ws.Range("K61").Value = 2935.55
$ws.Range("L61").Value = 3976.25
$ws.Range("M61").Value = -2723.55
$ws.Range("N61").Value = -4400.25

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0

$ws.Range("H132").Value = 1939.4
$ws.Range("I132").Value = 1713.4286
$ws.Range("J132").Value = 2466.6667
$ws.Range("K132").Value = 5140.2858
$ws.Range("L132").Value = 7400.000100000001
$ws.Range("M132").Value = -2610.2858
$ws.Range("N132").Value = -12460.0001

$ws.Range("H136").Value = 3109
$ws.Range("I136").Value = 2935.55
$ws.Range("J136").Value = 3976.25
$ws.Range("K136").Value = 8806.650000000001
$ws.Range("L136").Value = 11928.75
$ws.Range("M136").Value = -6256.650000000001
$ws.Range("N136").Value = -17028.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5009.4443
$ws.Range("I20").Value = 6565.8335
$ws.Range("K20").Value = 6565.8335
$ws.Range("M20").Value = -6318.8335

$ws.Range("H86").Value = 4248
$ws.Range("I86").Value = 3821.5
$ws.Range("J86").Value = 6807
$ws.Range("K86").Value = 3821.5
$ws.Range("L86").Value = 6807
$ws.Range("M86").Value = -2698.5
$ws.Range("N86").Value = -9053

$ws.Range("H89").Value = 4248
$ws.Range("I89").Value = 3821.5
$ws.Range("J89").Value = 6807
$ws.Range("K89").Value = 19107.5
$ws.Range("L89").Value = 34035
$ws.Range("M89").Value = -13491.5
$ws.Range("N89").Value = -45267

$ws.Range("H99").Value = 2628.2144
$ws.Range("I99").Value = 3576.5454
$ws.Range("J99").Value = 2014.5883
$ws.Range("K99").Value = 3576.5454
$ws.Range("L99").Value = 2014.5883
$ws.Range("M99").Value = -2078.5454
$ws.Range("N99").Value = -5010.588299999999

$ws.Range("H134").Value = 2349.125
$ws.Range("I134").Value = 2006.7693
$ws.Range("K134").Value = 6020.3079
$ws.Range("M134").Value = -3485.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 166.72
$ws.Range("I7").Value = 121.40909
$ws.Range("J7").Value = 499
$ws.Range("K7").Value = 121.40909
$ws.Range("L7").Value = 499
$ws.Range("M7").Value = -8.409090000000006
$ws.Range("N7").Value = -725

$ws.Range("H31").Value = 4239.1904
$ws.Range("J31").Value = 5374.933
$ws.Range("L31").Value = 5374.933
$ws.Range("N31").Value = -5964.933

$ws.Range("H34").Value = 4239.1904
$ws.Range("J34").Value = 5374.933
$ws.Range("L34").Value = 5374.933
$ws.Range("N34").Value = -5778.933

$ws.Range("H58").Value = 3443.9443
$ws.Range("I58").Value = 1955.2
$ws.Range("K58").Value = 1955.2
$ws.Range("M58").Value = -1752.2

$ws.Range("H60").Value = 16540.285
$ws.Range("I60").Value = 16540.285
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 16540.285
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -16029.285

$ws.Range("H99").Value = 14218.782
$ws.Range("J99").Value = 16644.076
$ws.Range("L99").Value = 16644.076
$ws.Range("N99").Value = -19640.076

$ws.Range("H105").Value = 1843.8
$ws.Range("I105").Value = 1269.6666
$ws.Range("J105").Value = 2705
$ws.Range("K105").Value = 1269.6666
$ws.Range("L105").Value = 2705
$ws.Range("M105").Value = 477.3334
$ws.Range("N105").Value = -6199

$ws.Range("H126").Value = 14218.782
$ws.Range("J126").Value = 16644.076
$ws.Range("L126").Value = 49932.228
$ws.Range("N126").Value = -54872.228

$ws.Range("H134").Value = 2375.4211
$ws.Range("I134").Value = 1942.4667
$ws.Range("J134").Value = 3999
$ws.Range("K134").Value = 5827.4001
$ws.Range("L134").Value = 11997
$ws.Range("M134").Value = -3292.4001
$ws.Range("N134").Value = -17067

$ws.Range("H136").Value = 3443.9443
$ws.Range("I136").Value = 1955.2
$ws.Range("K136").Value = 5865.6
$ws.Range("M136").Value = -3315.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 844
$ws.Range("J103").Value = 983.6667
$ws.Range("L103").Value = 2951.0001
$ws.Range("N103").Value = -4709.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2499.5
$ws.Range("I80").Value = 2249.5
$ws.Range("J80").Value = 2624.5
$ws.Range("K80").Value = 2249.5
$ws.Range("L80").Value = 2624.5
$ws.Range("M80").Value = -1251.5
$ws.Range("N80").Value = -4620.5

$ws.Range("H83").Value = 2499.5
$ws.Range("I83").Value = 2249.5
$ws.Range("J83").Value = 2624.5
$ws.Range("K83").Value = 11247.5
$ws.Range("L83").Value = 13122.5
$ws.Range("M83").Value = -6255.5
$ws.Range("N83").Value = -23106.5

$ws.Range("H97").Value = 888.6667
$ws.Range("I97").Value = 899.8
$ws.Range("K97").Value = 899.8
$ws.Range("M97").Value = -403.8

$ws.Range("H132").Value = 2841.9644
$ws.Range("I132").Value = 1890.5
$ws.Range("J132").Value = 3793.4285
$ws.Range("K132").Value = 5671.5
$ws.Range("L132").Value = 11380.2855
$ws.Range("M132").Value = -3141.5
$ws.Range("N132").Value = -16440.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2348.4666
$ws.Range("J82").Value = 1574.75
$ws.Range("L82").Value = 1574.75
$ws.Range("N82").Value = -2296.75

$ws.Range("H85").Value = 2348.4666
$ws.Range("J85").Value = 1574.75
$ws.Range("L85").Value = 1574.75
$ws.Range("N85").Value = -4070.75

$ws.Range("H93").Value = 5001.5
$ws.Range("I93").Value = 5001.5
$ws.Range("K93").Value = 5001.5
$ws.Range("M93").Value = -3753.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3133.2856
$ws.Range("I122").Value = 3143.5386
$ws.Range("K122").Value = 9430.6158
$ws.Range("M122").Value = -6980.6158

$ws.Range("H132").Value = 1446.4
$ws.Range("I132").Value = 1360.5385
$ws.Range("K132").Value = 4081.6155
$ws.Range("M132").Value = -1551.6155

$ws.Range("H136").Value = 3247.889
$ws.Range("I136").Value = 3247.889
$ws.Range("K136").Value = 9743.667000000001
$ws.Range("M136").Value = -7193.667000000001
